# Relabel the "Measure" column (B) values to include units, and widen the
# column to fit the new, longer text.
#
# Old labels -> New labels
#   Trial 1 Time      -> Trial 1 Time (seconds)
#   Trial 1 Exertion  -> Trial 1 Exertion (RPE)
#   Trial 2 Time      -> Trial 2 Time (seconds)
#   Trial 2 Exertion  -> Trial 2 Exertion (RPE)
#
# Each label appears once per runner (5 runners, rows 2-21).
# Updating every occurrence of a given old label lets the rewritten
# shared-strings table drop the now-unused old entry and append the new
# one, matching how Excel itself maintains the shared string table.
#
# The update order below (Exertion labels first, then Time labels) controls
# the order in which the new strings are appended to the shared string
# table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trial 1 Exertion -> Trial 1 Exertion (RPE)  (rows 3, 7, 11, 15, 19)
$ws.Cells.Item(3, 2).Value = "Trial 1 Exertion (RPE)"
$ws.Cells.Item(7, 2).Value = "Trial 1 Exertion (RPE)"
$ws.Cells.Item(11, 2).Value = "Trial 1 Exertion (RPE)"
$ws.Cells.Item(15, 2).Value = "Trial 1 Exertion (RPE)"
$ws.Cells.Item(19, 2).Value = "Trial 1 Exertion (RPE)"

# Trial 2 Exertion -> Trial 2 Exertion (RPE)  (rows 5, 9, 13, 17, 21)
$ws.Cells.Item(5, 2).Value = "Trial 2 Exertion (RPE)"
$ws.Cells.Item(9, 2).Value = "Trial 2 Exertion (RPE)"
$ws.Cells.Item(13, 2).Value = "Trial 2 Exertion (RPE)"
$ws.Cells.Item(17, 2).Value = "Trial 2 Exertion (RPE)"
$ws.Cells.Item(21, 2).Value = "Trial 2 Exertion (RPE)"

# Trial 1 Time -> Trial 1 Time (seconds)  (rows 2, 6, 10, 14, 18)
$ws.Cells.Item(2, 2).Value = "Trial 1 Time (seconds)"
$ws.Cells.Item(6, 2).Value = "Trial 1 Time (seconds)"
$ws.Cells.Item(10, 2).Value = "Trial 1 Time (seconds)"
$ws.Cells.Item(14, 2).Value = "Trial 1 Time (seconds)"
$ws.Cells.Item(18, 2).Value = "Trial 1 Time (seconds)"

# Trial 2 Time -> Trial 2 Time (seconds)  (rows 4, 8, 12, 16, 20)
$ws.Cells.Item(4, 2).Value = "Trial 2 Time (seconds)"
$ws.Cells.Item(8, 2).Value = "Trial 2 Time (seconds)"
$ws.Cells.Item(12, 2).Value = "Trial 2 Time (seconds)"
$ws.Cells.Item(16, 2).Value = "Trial 2 Time (seconds)"
$ws.Cells.Item(20, 2).Value = "Trial 2 Time (seconds)"

# Widen column B so the longer labels are fully visible (was ~13.45 chars,
# now roughly twice as wide).
$ws.Columns.Item(2).ColumnWidth = 27.8
